# Actualizacion automatica 2025-11-03 08:30:05
#
# Monthly rollover for "RIOS CARRION ANGEL BENIGNO":
#  - "VENTAS POR GRUPO": the per-product-group figures that had been
#    accumulated for the asesor are cleared back to 0 (period closed out),
#    and the "X de 24" fulfilled-client counters on the totals row drop to
#    match.
#  - "VENTA MENSUAL": the rolling 4-month window shifts forward by one
#    month (julio/agosto/septiembre/octubre -> agosto/septiembre/
#    octubre/noviembre); each client's figures shift left into the new
#    column layout and the freshly opened month starts at 0.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

# --- Sheet "VENTAS POR GRUPO": zero-out previously reported nonzero cells ---
$ws1.Range("G4").Value = 0
$ws1.Range("H4").Value = 0
$ws1.Range("M4").Value = 0
$ws1.Range("N4").Value = 0
$ws1.Range("P4").Value = 0
$ws1.Range("H5").Value = 0
$ws1.Range("P5").Value = 0
$ws1.Range("M8").Value = 0
$ws1.Range("D10").Value = 0
$ws1.Range("L10").Value = 0
$ws1.Range("M12").Value = 0
$ws1.Range("L17").Value = 0
$ws1.Range("M21").Value = 0
$ws1.Range("P21").Value = 0
$ws1.Range("L22").Value = 0
$ws1.Range("M22").Value = 0
$ws1.Range("M25").Value = 0

# Totals row: recount of "<n> de 24" clients-with-sales per column
$ws1.Range("D26").Value = "0 de 24"
$ws1.Range("G26").Value = "0 de 24"
$ws1.Range("H26").Value = "0 de 24"
$ws1.Range("L26").Value = "0 de 24"
$ws1.Range("M26").Value = "0 de 24"
$ws1.Range("N26").Value = "0 de 24"
$ws1.Range("P26").Value = "0 de 24"

# --- Sheet "VENTA MENSUAL": shift month columns forward by one (headers + data) ---
$ws2.Range("C1").Value = "agosto"
$ws2.Range("D1").Value = "septiembre"
$ws2.Range("E1").Value = "octubre"
$ws2.Range("F1").Value = "noviembre"

$ws2.Range("C3").Value = 0
$ws2.Range("D3").Value = -3519.22
$ws2.Range("E3").Value = 0

$ws2.Range("D4").Value = 306.24
$ws2.Range("E4").Value = 2548.88
$ws2.Range("F4").Value = 0

$ws2.Range("C5").Value = 0
$ws2.Range("E5").Value = 1362.43
$ws2.Range("F5").Value = 0

$ws2.Range("C6").Value = 4392.44
$ws2.Range("D6").Value = 1265.01
$ws2.Range("E6").Value = 0

$ws2.Range("C8").Value = 0
$ws2.Range("D8").Value = 2411.69
$ws2.Range("E8").Value = 663.55
$ws2.Range("F8").Value = 0

$ws2.Range("C9").Value = 2161.81
$ws2.Range("D9").Value = 142.56
$ws2.Range("E9").Value = 0

$ws2.Range("C10").Value = 366.34
$ws2.Range("D10").Value = 549.5
$ws2.Range("E10").Value = 4220.84
$ws2.Range("F10").Value = 0

$ws2.Range("C11").Value = 253.44
$ws2.Range("D11").Value = -3989.12
$ws2.Range("E11").Value = 0

$ws2.Range("C12").Value = 1090.97
$ws2.Range("D12").Value = 17655.41
$ws2.Range("E12").Value = 6935.82
$ws2.Range("F12").Value = 0

$ws2.Range("C13").Value = 0

$ws2.Range("D16").Value = 829.44
$ws2.Range("E16").Value = 0

$ws2.Range("C17").Value = 0
$ws2.Range("E17").Value = 3162.93
$ws2.Range("F17").Value = 0

$ws2.Range("D20").Value = 3252.41
$ws2.Range("E20").Value = 0

$ws2.Range("C21").Value = -33.7
$ws2.Range("D21").Value = 1632.93
$ws2.Range("E21").Value = 1687.32
$ws2.Range("F21").Value = 0

$ws2.Range("C22").Value = 274.75
$ws2.Range("D22").Value = -347.92
$ws2.Range("E22").Value = 14679.01
$ws2.Range("F22").Value = 0

$ws2.Range("C25").Value = -81.41
$ws2.Range("D25").Value = 3690.09
$ws2.Range("E25").Value = 6348.54
$ws2.Range("F25").Value = 0

# Totals row (column sums)
$ws2.Range("C26").Value = 8424.64
$ws2.Range("D26").Value = 23879.02
$ws2.Range("E26").Value = 41609.32
$ws2.Range("F26").Value = 0

# Column widths re-auto-fit themselves after the data shift (Excel's
# ColumnWidth property is expressed in characters and is offset by the
# engine's default-font padding, so subtract that 5/6-character pad to
# land on the target stored widths of 13 / 16 / 14 / 15).
$ws2.Columns.Item(3).ColumnWidth = 12.166666666666666
$ws2.Columns.Item(4).ColumnWidth = 15.166666666666666
$ws2.Columns.Item(5).ColumnWidth = 13.166666666666666
$ws2.Columns.Item(6).ColumnWidth = 14.166666666666666
